$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 3 new rows at 13, pushing the old row 13 (linkflow) down to row 16 ---
$ws.Rows("13:15").Insert()

# --- Column F: switch the "real multiplier" formula from C/25 to E*D for rows 3-12 ---
$ws.Range("F3").Formula  = "=E3*D3"
$ws.Range("F4").Formula  = "=E4*D4"
$ws.Range("F5").Formula  = "=E5*D5"
$ws.Range("F6").Formula  = "=E6*D6"
$ws.Range("F7").Formula  = "=E7*D7"
$ws.Range("F8").Formula  = "=E8*D8"
$ws.Range("F9").Formula  = "=E9*D9"
$ws.Range("F10").Formula = "=E10*D10"
$ws.Range("F11").Formula = "=E11*D11"
$ws.Range("F12").Formula = "=E12*D12"

# --- HorseCoding (row 8): price updated manually, 18 -> 22.5 (now a literal, not a formula) ---
$ws.Range("C8").Value = 22.5

# --- FoxCode (row 12): price updated manually, 112.5 -> 150 (now a literal, not a formula) ---
$ws.Range("C12").Value = 150

# --- FoxCode tiered plans: new rows 13-15 ---
$ws.Range("A13").Value = "[FoxCode(狐狸)](https://foxcode.rjj.cc/auth/register?aff=6W5J31UI)"
$ws.Range("B13").Value = "富可敌国"
$ws.Range("C13").Value = 150
$ws.Range("D13").Value = 0.35
$ws.Range("E13").Value = 6
$ws.Range("F13").Formula = "=E13*D13"
$ws.Range("G13").Value = 35
$ws.Range("H13").Value = "无"
$ws.Range("I13").Value = "无"

$ws.Range("A14").Value = "[FoxCode(狐狸)](https://foxcode.rjj.cc/auth/register?aff=6W5J31UI)"
$ws.Range("B14").Value = "富可敌国"
$ws.Range("C14").Value = 150
$ws.Range("D14").Value = 0.27
$ws.Range("E14").Value = 6
$ws.Range("F14").Formula = "=E14*D14"
$ws.Range("G14").Value = 135
$ws.Range("H14").Value = "无"
$ws.Range("I14").Value = "无"

$ws.Range("A15").Value = "[FoxCode(狐狸)](https://foxcode.rjj.cc/auth/register?aff=6W5J31UI)"
$ws.Range("B15").Value = "富可敌国"
$ws.Range("C15").Value = 150
$ws.Range("D15").Value = 0.23
$ws.Range("E15").Value = 6
$ws.Range("F15").Formula = "=E15*D15"
$ws.Range("G15").Value = 468
$ws.Range("H15").Value = "无"
$ws.Range("I15").Value = "无"

# --- Row 16 (old row 13, shifted down by the insert): update its F formula to match the new pattern ---
$ws.Range("F16").Formula = "=E16*D16"

# --- Sheet view cosmetics: zoom + selection ---
$ws.Range("F15").Select()
$excel.ActiveWindow.Zoom = 150
